$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Canada")
Write-Output $ws.Name
